# Fixed update to excel issue
#
# 1. Rename header cells:
#    - "Weekly Quantity"!B1  "Requested quantity" -> "Weekly_PO_Qty"
#    - "Monthly Trend"!B1    "Requested quantity" -> "Monthly_PO_Qty"
# 2. Add a new "PO Forecast" worksheet (after "Monthly Trend") containing
#    a Prophet-style forecast table: ds, PO_Forecast, yhat_lower, yhat_upper

$wb = $excel.ActiveWorkbook

# --- 1. Rename existing headers ---------------------------------------
$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet -
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

# Header row
$forecast.Cells.Item(1, 1).Value = "ds"
$forecast.Cells.Item(1, 2).Value = "PO_Forecast"
$forecast.Cells.Item(1, 3).Value = "yhat_lower"
$forecast.Cells.Item(1, 4).Value = "yhat_upper"

$headerRange = $forecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous
$headerRange.Borders.Weight = 2           # xlThin

# Data rows: ds (weekly date serials), PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45137.99999999999, 7, 4.482587678229596, 9.027211824315206),
    @(45151.99999999999, 6, 3.344174294339976, 8.050919866023959),
    @(45172.99999999999, 4, 1.673604461542263, 6.096201250260933),
    @(45186.99999999999, 3, 0.4775874275924933, 5.045401035077584),
    @(45193.99999999999, 2, -0.0973995051111579, 4.55939883598687),
    @(45200.99999999999, 2, -0.6666137592989833, 3.829218949746277),
    @(45207.99999999999, 1, -1.231056601802954, 3.330307192156785),
    @(45214.99999999999, 1, -1.579561432530414, 2.821672044553867),
    @(45221.99999999999, 0, -2.390998330916422, 2.199268123350321),
    @(45228.99999999999, 0, -2.780049271951468, 1.639181243994527),
    @(45235.99999999999, 0, -3.432283297365349, 0.9957509155046844),
    @(45242.99999999999, 0, -3.945713042406078, 0.5552170458640834)
)

$row = 2
foreach ($r in $data) {
    $dateCell = $forecast.Cells.Item($row, 1)
    $dateCell.Value = $r[0]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $forecast.Cells.Item($row, 2).Value = $r[1]
    $forecast.Cells.Item($row, 3).Value = $r[2]
    $forecast.Cells.Item($row, 4).Value = $r[3]

    $row++
}

# Restore the originally active sheet/selection so we don't disturb the
# workbook's active-tab state just by virtue of having added a new sheet.
[void]$weekly.Activate()
[void]$weekly.Range("A1").Select()

Write-Host "Sheets now: $([string]::Join(', ', ($wb.Worksheets | ForEach-Object { $_.Name })))"
